# Apply updated crypto price/volume figures (text cells; values are strings, not numbers).
# A leading apostrophe forces Excel to keep the assigned text as-is (no numeric/
# percentage auto-conversion), and resetting Style to "Normal" avoids leaving a
# stray quote-prefix cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "28.310.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "  +0.72%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + "1.857.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "  -0.47%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + "  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + "330.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + "  -1.91%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'" + "0.4537"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + "  -3.48%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'" + "0.3907"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + "  +0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + "47.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + "  +1.78%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + "0.07777"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + "  -2.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + "0.9788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + "  -1.04%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + "21.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "  -1.46%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + "1.843.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "  +0.56%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + "5.762"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + "  -3.47%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'" + "6.936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + "  -4.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + "  -0.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + "87.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + "  -4.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'" + "0.06516"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + "  -1.53%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + "0.00001014"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + "  -2.80%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + "16.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + "  -3.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + "  +0.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + "28.287.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + "  +0.62%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + "5.263"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + "  -2.80%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + "10.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + "  -3.72%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + "2.249"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + "  -1.62%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + "2.069.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + "  +0.63%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + "156.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + "  -1.95%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'" + "19.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + "  -2.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'" + "  -4.25%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + "5.250"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + "  -4.34%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + "116.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + "  -2.82%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'" + "0.09243"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + "  -2.66%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + "0.9327"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + "  -4.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + "3.598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + "  +0.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + "1.372"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + "  +0.94%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + "5.170"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + "  -2.90%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + "0.05990"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + "  -1.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + "0.02181"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + "  -3.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'" + "  -2.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + "1.162"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "  -0.72%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + "  -0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + "0.5633"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + "  -5.43%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + "9.933"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + "  -3.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'" + "0.1780"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + "  -5.04%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + "1.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + "  -1.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + "2.315"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + "  +24.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + "0.07193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + "  +4.63%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + "0.5347"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + "  -4.42%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'" + "  -3.83%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'" + "1.859"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + "  -5.66%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + "109.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + "  -2.38%  "
$ws.Range("E51").Style = "Normal"
